$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text in E8 ("Good Morning" -> "GIT UPDATE")
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the user's active selection on that cell
$ws.Range("E8").Select() | Out-Null
